# Normalizes the "rol_participante" (column T) values from upper-case /
# snake_case codes to proper Title Case text, and fixes a set of
# zero-value "edad" (column V) placeholders to 30 (data-cleaning pass -
# "Readme y Dashboard terminados, Proyecto finalizado").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build accented replacement strings via char codes so the source file
# (which must stay plain ASCII) round-trips correctly through the
# interpreter's string literal handling.
$nTilde = [char]0x00F1
$NTilde = [char]0x00D1

$pasajeroKey = "PASAJERO_ACOMPA" + $NTilde + "ANTE"
$pasajeroVal = "Pasajero Acompa" + $nTilde + "ante"

$roleMap = @{
    "CONDUCTOR" = "Conductor"
    "CICLISTA" = "Ciclista"
}
$roleMap[$pasajeroKey] = $pasajeroVal

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $roleCell = $ws.Cells.Item($r, 20)   # column T = rol_participante
    $roleVal = $roleCell.Text
    if ($roleMap.ContainsKey($roleVal)) {
        $roleCell.Value = $roleMap[$roleVal]
    }

    $ageCell = $ws.Cells.Item($r, 22)    # column V = edad
    $ageVal = $ageCell.Text
    if ($ageVal -eq "0") {
        $ageCell.Value = 30
    }
}

Write-Output "done"
